$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the LocationID column entirely (column A), shifting ProductID /
# MaterialPrice left by one column.
$ws.Columns.Item(1).Delete()

# Match the post-edit selection seen in the source workbook.
$ws.Range("O12").Select()
